# Update header labels so the first row can be used automatically as a
# header when the table is loaded into Power BI.
$wb = $excel.ActiveWorkbook

# Sheets 1, 2, 3 and 5 use a simple "Ano <year>" label pattern.
$anoSheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Emissoes Totais (MtCO2eq)"
)

foreach ($name in $anoSheets) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("B1").Value = "Ano " + $ws.Range("B1").Text
    $ws.Range("C1").Value = "Ano " + $ws.Range("C1").Text
    $ws.Range("D1").Value = "Ano " + $ws.Range("D1").Text
    $ws.Range("E1").Value = "Ano " + $ws.Range("E1").Text
}

# Sheet 4 uses an "Intervalo <range>" label pattern.
$ws4 = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
$ws4.Range("B1").Value = "Intervalo " + $ws4.Range("B1").Text
$ws4.Range("C1").Value = "Intervalo " + $ws4.Range("C1").Text
$ws4.Range("D1").Value = "Intervalo " + $ws4.Range("D1").Text
$ws4.Range("E1").Value = "Intervalo " + $ws4.Range("E1").Text

# Sheet 6 only has a single year column (B1).
$ws6 = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$ws6.Range("B1").Value = "Ano " + $ws6.Range("B1").Text
